# Commit: "Added converter (fix int cell issue)"
# A new "Column with int" column (F) is added to Sheet1 with integer values,
# and the active sheet/selection state is updated to reflect Sheet1 being
# the active tab with G5 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New column F: header + two integer values
$ws1.Range("F1").Value = "Column with int"
$ws1.Range("F2").Value = 10
$ws1.Range("F3").Value = 12

# Sheet1 becomes the active sheet with G5 selected (also clears the
# previous "Test" sheet's tab-selected/active state)
$ws1.Activate()
$ws1.Range("G5").Select() | Out-Null
